# Applies the authored edit to "Capacity Supply Curve.xlsx":
#   1. About!C1 date bumped from 2/8/2024 (45330) to 3/11/2024 (45362)
#   2. CSC-CSCSoCECBiaSY: the "Share of Cost Effective Capacity Built in a
#      Single Year" table rows are set to 1 (was 0.33 / 0.25) for every
#      plant type that had existing capacity (rows 2-15 and 18-25); the two
#      rows with no existing capacity (16-17) are left at 0, unchanged.
#   3. The CSC-CSCSoCECBiaSY sheet's on-screen selection is moved to
#      B18:AE25, matching where the author was last working.

$wb = $excel.ActiveWorkbook

# --- 1. About sheet: bump the "last updated" date -------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45362

# --- 2. CSC-CSCSoCECBiaSY: set calibrated share values to 1 ----------------
$wsShare = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

# Rows with existing capacity (their share is being set to 1 = 100%).
$rowsToOne = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,18,19,20,21,22,23,24,25)
foreach ($r in $rowsToOne) {
    $wsShare.Range("B" + $r + ":AE" + $r).Value = 1
}

# --- 3. Update the active selection on CSC-CSCSoCECBiaSY -------------------
$wsShare.Activate()
$wsShare.Range("B18:AE25").Select()
